# Updated cryptos list on Sat Feb 18 12:15:00 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for the cryptos table, and
# swaps the Chainlink / ShibaInu rows (15 <-> 16) to reflect their new rank
# order. Price strings look numeric (e.g. "1.000", "0.3936") but must stay
# as literal text like the source sheet, so each Price cell is forced to a
# text number-format before the write and reset back to the Normal style
# afterwards (Excel would otherwise silently coerce them to numbers and
# drop the significant trailing/leading zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Rows 15 and 16 swap rank order: Chainlink moves up to rank 15 (row 15),
# ShibaInu moves down to rank 16 (row 16). Both get refreshed price/volume data.
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-PriceText "D15" "7.790"
$ws.Range("E15").Value = "  +6.67%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-PriceText "D16" "0.00001314"
$ws.Range("E16").Value = "  -0.37%  "

# Refresh prices/volumes for the remaining rows.
Set-PriceText "D2" "24.560.58"
$ws.Range("E2").Value = "  +3.18%  "

Set-PriceText "D3" "1.692.88"
$ws.Range("E3").Value = "  +1.53%  "

Set-PriceText "D4" "1.000"
$ws.Range("E4").Value = "  +0.04%  "

Set-PriceText "D5" "315.51"
$ws.Range("E5").Value = "  +1.58%  "

Set-PriceText "D6" "1.000"
$ws.Range("E6").Value = "  +0.06%  "

Set-PriceText "D7" "0.3936"
$ws.Range("E7").Value = "  +1.36%  "

Set-PriceText "D8" "0.3999"
$ws.Range("E8").Value = "  +0.90%  "

Set-PriceText "D9" "1.525"
$ws.Range("E9").Value = "  +4.45%  "

Set-PriceText "D10" "1.000"
$ws.Range("E10").Value = "  +0.03%  "

Set-PriceText "D11" "52.63"
$ws.Range("E11").Value = "  +2.99%  "

Set-PriceText "D12" "0.08742"
$ws.Range("E12").Value = "  +0.70%  "

Set-PriceText "D13" "7.226"
$ws.Range("E13").Value = "  +6.82%  "

Set-PriceText "D14" "23.12"
$ws.Range("E14").Value = "  +2.12%  "

Set-PriceText "D17" "1.698.36"
$ws.Range("E17").Value = "  +1.77%  "

Set-PriceText "D18" "99.45"
$ws.Range("E18").Value = "  -0.12%  "

Set-PriceText "D19" "0.07079"
$ws.Range("E19").Value = "  +3.79%  "

Set-PriceText "D20" "19.61"
$ws.Range("E20").Value = "  +2.39%  "

Set-PriceText "D21" "6.884"
$ws.Range("E21").Value = "  +3.44%  "

$ws.Range("E22").Value = "  +0.19%  "

Set-PriceText "D23" "14.04"
$ws.Range("E23").Value = "  +1.09%  "

Set-PriceText "D24" "24.555.97"
$ws.Range("E24").Value = "  +3.20%  "

Set-PriceText "D25" "3.128"
$ws.Range("E25").Value = "  +7.00%  "

Set-PriceText "D26" "2.347"
$ws.Range("E26").Value = "  +1.50%  "

Set-PriceText "D27" "22.28"
$ws.Range("E27").Value = "  +2.33%  "

Set-PriceText "D28" "161.27"
$ws.Range("E28").Value = "  +1.29%  "

Set-PriceText "D29" "5.215"
$ws.Range("E29").Value = "  +1.38%  "

Set-PriceText "D30" "134.64"
$ws.Range("E30").Value = "  +3.67%  "

Set-PriceText "D31" "7.614"
$ws.Range("E31").Value = "  +12.15%  "

Set-PriceText "D32" "1.878.25"
$ws.Range("E32").Value = "  +1.45%  "

$ws.Range("E33").Value = "  -3.30%  "

Set-PriceText "D34" "7.361"
$ws.Range("E34").Value = "  +11.17%  "

Set-PriceText "D35" "0.08537"
$ws.Range("E35").Value = "  +0.10%  "

Set-PriceText "D36" "11.32"
$ws.Range("E36").Value = "  +7.96%  "

Set-PriceText "D37" "1.935"
$ws.Range("E37").Value = "  -0.65%  "

Set-PriceText "D38" "0.2711"

Set-PriceText "D39" "14.44"
$ws.Range("E39").Value = "  -0.52%  "

Set-PriceText "D40" "0.02741"
$ws.Range("E40").Value = "  +8.87%  "

Set-PriceText "D41" "0.09047"
$ws.Range("E41").Value = "  +2.61%  "

Set-PriceText "D42" "1.479"
$ws.Range("E42").Value = "  +0.93%  "

Set-PriceText "D43" "0.7663"
$ws.Range("E43").Value = "  +0.89%  "

Set-PriceText "D44" "0.7161"
$ws.Range("E44").Value = "  +1.67%  "

Set-PriceText "D45" "15.34"
$ws.Range("E45").Value = "  +2.48%  "

Set-PriceText "D46" "2.520"
$ws.Range("E46").Value = "  +3.22%  "

Set-PriceText "D47" "4.206"
$ws.Range("E47").Value = "  +2.46%  "

Set-PriceText "D48" "1.0000"
$ws.Range("E48").Value = "  +0.06%  "

Set-PriceText "D49" "1.336"
$ws.Range("E49").Value = "  +9.70%  "

Set-PriceText "D50" "141.16"
$ws.Range("E50").Value = "  +1.12%  "

Set-PriceText "D51" "0.07987"
$ws.Range("E51").Value = "  +2.55%  "
